$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RMSE value in A1
$ws.Range("A1").Value = 45.715882145518449
